# #12 SAI everytime added on slides
#
# 1) Insert a new "SAI everytime" slide at position 10 (Title and Content
#    layout - same layout used by the neighbouring "Add"/"Solutions" slides).
# 2) Refresh the cached "datetimeFigureOut" footer field (3/7/2020 -> 3/8/2020)
#    on the slide master and on every slide layout.

$p = $ppt.ActivePresentation

# --- 1) New slide -----------------------------------------------------
$s = $p.Slides.Add(10, 2)

$title = $s.Shapes.Item(1)
$title.Name = "Titre 1"
$title.TextFrame.TextRange.Text = "SAI everytime"
$title.TextFrame.TextRange.LanguageID = "fr-FR"

$body = $s.Shapes.Item(2)
$body.Name = "Espace réservé du contenu 2"
$bodyText = "SAI need to be launch everytime the computer is up, it needs to run as a daemon that use the mic and let the mic for another application if its need" + [char]13 + "A schema to explain"
$body.TextFrame.TextRange.Text = $bodyText
$body.TextFrame.TextRange.LanguageID = "fr-FR"

# --- 2) Footer date field refresh -------------------------------------
$newDate = "3/8/2020"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
